# Auto-generated edit script applying the diff changes
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("G4").Value = 60
$ws.Range("F5").Value = 62
$ws.Range("G5").Value = 29.9
$ws.Range("G7").Value = '不可售'
$ws.Range("G8").Value = 60
$ws.Range("G9").Value = '不可售'
$ws.Range("F10").Value = 392
$ws.Range("F13").Value = 325
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 592
$ws.Range("F19").Value = 1485
$ws.Range("F20").Value = 5782
$ws.Range("F22").Value = 1629
$ws.Range("F24").Value = 77
$ws.Range("F25").Value = 37
$ws.Range("F26").Value = 5442
$ws.Range("F27").Value = 5442
$ws.Range("F30").Value = 1559
$ws.Range("F31").Value = 296
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 70
$ws.Range("F34").Value = 1089
$ws.Range("F36").Value = 135
$ws.Range("F37").Value = 9

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 177
$ws.Range("F8").Value = 254

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 9432
$ws.Range("F4").Value = 2167
$ws.Range("F5").Value = 279

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 9432
$ws.Range("F4").Value = 2167
$ws.Range("G7").Value = 60
$ws.Range("F8").Value = 62
$ws.Range("G8").Value = 29.9
$ws.Range("C10").Value = '杭州·萌忧·原崩铁同人only'
$ws.Range("D10").Value = '康候圣街99号 顺丰创新中心'
$ws.Range("E10").Value = '2024.08.24 10:30-08.24 17:00'
$ws.Range("F10").Value = 1082
$ws.Range("G10").Value = 60
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=87293'
$ws.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202406/rQFz5smR1717475284585.jpeg'
$ws.Range("C11").Value = '杭州·首届Fun-X动漫嘉年华【免费入场】'
$ws.Range("D11").Value = '文三路 玩美的一天沉浸式生活街区'
$ws.Range("E11").Value = '2024.08.24 09:30-08.24 17:30'
$ws.Range("F11").Value = 392
$ws.Range("G11").Value = 60
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=89710'
$ws.Range("I11").Value = '//i1.hdslb.com/bfs/openplatform/202408/2FEfmcIE1722914643220.jpeg'
$ws.Range("B12").Value = '2024-08-25'
$ws.Range("C12").Value = '杭州·代号鸢同人only-春风夜夜电音夜场'
$ws.Range("D12").Value = '三墩镇蒋墩路375号(华策中心A座) 经典时代Live House'
$ws.Range("E12").Value = '2024.08.25 17:00-08.26 00:00'
$ws.Range("F12").Value = 438
$ws.Range("G12").Value = 158
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=90025'
$ws.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202408/dMGmIqkv1724228307637.jpeg'
$ws.Range("B13").Value = '2024-08-31'
$ws.Range("C13").Value = '杭州·初音未来17周年生日派对 & 音链视窗同人共创only'
$ws.Range("D13").Value = '金惠路1128号西区 杭州金迪大酒店'
$ws.Range("E13").Value = '2024.08.31 12:00-08.31 20:00'
$ws.Range("F13").Value = 325
$ws.Range("G13").Value = 39
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=90372'
$ws.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202408/j61I0I7n1722925055083.jpeg'
$ws.Range("C14").Value = '杭州·音乐番+only'
$ws.Range("D14").Value = '康候圣街99号 顺丰创新中心'
$ws.Range("E14").Value = '2024.08.31 10:00-08.31 18:00'
$ws.Range("F14").Value = 371
$ws.Range("G14").Value = 78
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=88899'
$ws.Range("I14").Value = '//i0.hdslb.com/bfs/openplatform/202408/lxW52TpT1724228135568.jpeg'
$ws.Range("B15").Value = '2024-09-06'
$ws.Range("C15").Value = '杭州·次元萌友会【免费展会】'
$ws.Range("D15").Value = '祥泰街398号 杭州万融城'
$ws.Range("E15").Value = '2024.09.06 10:00-09.08 21:00'
$ws.Range("F15").Value = 53
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=90896'
$ws.Range("I15").Value = '//i0.hdslb.com/bfs/openplatform/202408/snpy8ATR1723793956830.png'
$ws.Range("B16").Value = '2024-09-07'
$ws.Range("C16").Value = '杭州·DNP01综合同人展X【昼夜星逐】泛VOCALOID专场'
$ws.Range("D16").Value = '观澜路钱江世纪公园d区1幢 杭州世纪雷迪森庄园酒店'
$ws.Range("E16").Value = '2024.09.07 12:00-09.07 20:00'
$ws.Range("F16").Value = 67
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=90587'
$ws.Range("I16").Value = '//i0.hdslb.com/bfs/openplatform/202408/x9rSjkDf1723208824749.jpeg'
$ws.Range("C17").Value = '杭州·《卡农》永恒经典名曲音乐会'
$ws.Range("D17").Value = '曙光路31号 浙江音乐厅'
$ws.Range("E17").Value = '2024.09.07 19:30-09.07 21:00'
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 100
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=85894'
$ws.Range("I17").Value = '//i2.hdslb.com/bfs/openplatform/202405/3jz9YpaW1716100738530.jpeg'
$ws.Range("C18").Value = '杭州·红楼梦·主题演绎国风音乐会《梦寻红楼》'
$ws.Range("D18").Value = '望梅路与汀兰路交叉口向南100米 杭州临平大剧院（原余杭大剧院）'
$ws.Range("E18").Value = '2024.09.07 15:00-09.07 16:30'
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=89257'
$ws.Range("I18").Value = '//i2.hdslb.com/bfs/openplatform/202407/tkm6AHo71720572975141.jpeg'
$ws.Range("B19").Value = '2024-09-15'
$ws.Range("C19").Value = '杭州·2024首届COMIC GALAXY次元盛典'
$ws.Range("D19").Value = '长江南路336号 白马湖国际会展中心'
$ws.Range("E19").Value = '2024.09.15 09:30-09.17 17:30'
$ws.Range("F19").Value = 592
$ws.Range("G19").Value = 88
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=90433'
$ws.Range("I19").Value = '//i0.hdslb.com/bfs/openplatform/202408/teoBMbzd1723019674766.png'
$ws.Range("F20").Value = 1485
$ws.Range("F21").Value = 5782
$ws.Range("F23").Value = 1629
$ws.Range("F29").Value = 5442
$ws.Range("F30").Value = 5442
$ws.Range("F33").Value = 1559
$ws.Range("F34").Value = 297
$ws.Range("F35").Value = 29
$ws.Range("F36").Value = 1089
$ws.Range("F38").Value = 135
